$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All the cells in column B that currently read "[88]" (the repeated
# placeholder value) should now read "[89]" -- except B88, whose own
# index is being repointed separately below.
$cellsToBump = @(
    "B4","B11","B16","B18","B21","B23","B25","B28","B31","B33",
    "B35","B38","B40","B44","B48","B52","B56","B59","B62","B64",
    "B66","B68","B70","B72","B74","B76","B78","B81","B83","B85"
)
foreach ($addr in $cellsToBump) {
    $ws.Range($addr).Value = "[89]"
}

# B87 moves from "[87]" to "[88]".
$ws.Range("B87").Value = "[88]"

# Append a new row 89 at the bottom of the table.
$ws.Range("A89").Value = 88.0
$ws.Range("B89").Value = "[89]"
